$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 274
$ws1.Range("F3").Value = 256
$ws1.Range("F8").Value = 2253
$ws1.Range("F9").Value = 1463
$ws1.Range("F10").Value = 40
$ws1.Range("F15").Value = 1417
$ws1.Range("F16").Value = 5407
$ws1.Range("F18").Value = 5409
$ws1.Range("F19").Value = 1985
$ws1.Range("F20").Value = 2947
$ws1.Range("F21").Value = 3373
$ws1.Range("F22").Value = 194
$ws1.Range("F23").Value = 1640
$ws1.Range("F24").Value = 22
$ws1.Range("F25").Value = 275
$ws1.Range("F26").Value = 852
$ws1.Range("F27").Value = 147
$ws1.Range("F28").Value = 10
$ws1.Range("F29").Value = 337
$ws1.Range("F30").Value = 1054
$ws1.Range("F31").Value = 2191
$ws1.Range("F33").Value = 133
$ws1.Range("F34").Value = 310
$ws1.Range("B35").Value = "'2024-08-17"
$ws1.Range("E35").Value = "2024.08.17 09:00-08.18 17:00"
$ws1.Range("F35").Value = 807
$ws1.Range("G35").Value = 75
$ws1.Range("I35").Value = "//i2.hdslb.com/bfs/openplatform/202405/TU8kiduQ1715238040248.jpeg"
$ws1.Range("F37").Value = 390
$ws1.Range("F38").Value = 457

# ---- Sheet 2: 演出 ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F12").Value = 26
$ws2.Range("F13").Value = 200
$ws2.Range("F17").Value = 17
$ws2.Range("F19").Value = 89
$ws2.Range("F21").Value = 52

# ---- Sheet 2: new row 22 ----
$ws2.Range("A22").Value = 21
$ws2.Range("B22").Value = "'2024-08-09"
$ws2.Range("C22").Value = "北京·燃爆DNA——日本动漫原声金曲超燃演唱会"
$ws2.Range("D22").Value = "复兴门内大街49号 民族宫大剧院"
$ws2.Range("E22").Value = "2024.08.09 19:30-08.09 21:30"
$ws2.Range("F22").Value = 0
$ws2.Range("G22").Value = 180
$ws2.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=85334"
$ws2.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202405/WpZshtXD1715052832157.jpeg"
$ws2.Range("A21").Copy()
$ws2.Range("A22").PasteSpecial(-4122)

# ---- Sheet 4: 全部类型 ----
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F7").Value = 274
$ws4.Range("F12").Value = 2253
$ws4.Range("F13").Value = 1463
$ws4.Range("F14").Value = 40
$ws4.Range("F17").Value = 26
$ws4.Range("F19").Value = 1417
$ws4.Range("F20").Value = 200
$ws4.Range("F24").Value = 5407
$ws4.Range("F26").Value = 5409
$ws4.Range("F27").Value = 1985
$ws4.Range("F28").Value = 2947
$ws4.Range("F29").Value = 3373
$ws4.Range("F30").Value = 17
$ws4.Range("F31").Value = 194
$ws4.Range("F33").Value = 89
$ws4.Range("F34").Value = 1640
$ws4.Range("F36").Value = 275
$ws4.Range("F37").Value = 852
$ws4.Range("F38").Value = 147
$ws4.Range("F39").Value = 10
$ws4.Range("F40").Value = 337
$ws4.Range("F41").Value = 52
$ws4.Range("F42").Value = 2191
$ws4.Range("F44").Value = 133
$ws4.Range("F45").Value = 310
$ws4.Range("B46").Value = "'2024-08-17"
$ws4.Range("E46").Value = "2024.08.17 09:00-08.18 17:00"
$ws4.Range("F46").Value = 807
$ws4.Range("G46").Value = 75
$ws4.Range("I46").Value = "//i2.hdslb.com/bfs/openplatform/202405/TU8kiduQ1715238040248.jpeg"
$ws4.Range("F48").Value = 390
$ws4.Range("F49").Value = 457
